$d = $word.ActiveDocument

# Original paragraph text is "Version 1." with character positions:
#   0123456 7 8 9
#   Version   1 .
#
# Target structure:
#   <w:r>Versi</w:r><w:r>on</w:r> ... " 2" ... (bookmark) ... <w:r>.</w:r>
#
# i.e. "Version" is split into two runs ("Versi" + "on"), the version
# number changes from 1 to 2, and the trailing period is moved to its own
# run placed after the _GoBack bookmark.

# Step 1: split "Version" into "Versi" + "on" as two separate runs.
# Using Range.InsertXML on the "on" sub-range merges the supplied run(s)
# into the surrounding paragraph while keeping the preceding "Versi" text
# as its own, separate run (no leftover empty <w:rPr/>).
$rOn = $d.Range(5, 7)
$onXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
         '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
         '<pkg:xmlData>' +
         '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
         '<w:body><w:p><w:r><w:t>on</w:t></w:r></w:p></w:body>' +
         '</w:document>' +
         '</pkg:xmlData></pkg:part></pkg:package>'
$rOn.InsertXML($onXml)

# Step 2: change the version digit "1" -> "2" (still at positions 8-9).
$rDigit = $d.Range(8, 9)
$rDigit.Text = "2"

# Step 3: remove the trailing period that currently sits before the
# bookmark (positions 9-10).
$rPeriod = $d.Range(9, 10)
$rPeriod.Text = ""

# Step 4: retype the period using the Selection object positioned right
# after the (now trailing) bookmark; this creates a brand-new run placed
# after the bookmark rather than merging into the run before it.
$sel = $word.Selection
$sel.SetRange(9, 9)
$sel.TypeText(".")
